$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list values (Price / Volume(1h)) - stored as text so the
# numeric/percent-looking strings (e.g. "277.00", "1.33%") are preserved
# literally instead of being auto-converted to numbers by Excel.
$updates = @{
    'D2' = '277.00'
    'E2' = '1.33%'
    'D3' = '27.19'
    'E3' = '1.48%'
    'D4' = '4.922'
    'E4' = '0.50%'
    'E5' = '1.42%'
    'D6' = '6.925'
    'E6' = '0.46%'
    'D7' = '1.249'
    'E7' = '-6.91%'
    'D8' = '0.8833'
    'E8' = '-0.06%'
    'E9' = '3.83%'
    'D10' = '0.05082'
    'E10' = '-0.37%'
    'D11' = '0.07524'
    'E11' = '1.63%'
    'E12' = '-8.20%'
    'D13' = '0.09012'
    'E13' = '-0.25%'
    'D14' = '0.001579'
    'E14' = '0.34%'
    'D15' = '0.0006423'
    'E15' = '2.03%'
    'D16' = '0.006014'
    'E16' = '-0.45%'
    'E17' = '-0.28%'
    'D18' = '3.316'
    'E18' = '-1.17%'
    'D21' = '0.1328'
    'E21' = '-0.43%'
    'D22' = '3.918'
    'E22' = '0.28%'
    'D23' = '0.04426'
    'E23' = '2.16%'
    'D24' = '0.001176'
    'E24' = '-0.11%'
    'D25' = '0.003867'
    'E25' = '5.83%'
    'E26' = '-0.06%'
    'E27' = '14.04%'
    'D40' = '0.04134'
    'E40' = '2.42%'
    'D41' = '0.006815'
    'E41' = '3.05%'
    'D42' = '0.1180'
    'E42' = '1.36%'
    'D43' = '0.002231'
    'E43' = '4.62%'
    'D44' = '0.01125'
    'E44' = '-8.35%'
    'D45' = '0.00005201'
    'E45' = '-2.60%'
    'D46' = '1.483'
    'E46' = '-37.06%'
    'D47' = '0.02025'
    'E47' = '-4.48%'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
